# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates across 8 sheets per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 288.5
$ws.Range("I2").Value = 77
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 77
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = 36
$ws.Range("N2").Value = -726
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("H12").Value = 150
$ws.Range("I12").Value = 133.33333
$ws.Range("K12").Value = 133.33333
$ws.Range("M12").Value = 36.66667000000001
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("H19").Value = 1405.25
$ws.Range("I19").Value = 1120
$ws.Range("J19").Value = 1534.909
$ws.Range("K19").Value = 1120
$ws.Range("L19").Value = 1534.909
$ws.Range("M19").Value = -945
$ws.Range("N19").Value = -1884.909
$ws.Range("H31").Value = 45000
$ws.Range("I31").Value = 50000
$ws.Range("J31").Value = 40000
$ws.Range("K31").Value = 150000
$ws.Range("L31").Value = 120000
$ws.Range("M31").Value = -149770
$ws.Range("N31").Value = -120460
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H61").Value = 354
$ws.Range("I61").Value = 305.8
$ws.Range("J61").Value = 474.5
$ws.Range("K61").Value = 917.4000000000001
$ws.Range("L61").Value = 1423.5
$ws.Range("M61").Value = -745.4000000000001
$ws.Range("N61").Value = -1767.5
$ws.Range("H69").Value = 9776.333000000001
$ws.Range("J69").Value = 10283.857
$ws.Range("L69").Value = 30851.571
$ws.Range("N69").Value = -32599.571
$ws.Range("H72").Value = 9776.333000000001
$ws.Range("J72").Value = 10283.857
$ws.Range("L72").Value = 92554.713
$ws.Range("N72").Value = -101290.713
$ws.Range("H86").Value = 2319.923
$ws.Range("I86").Value = 1344.0667
$ws.Range("K86").Value = 1344.0667
$ws.Range("M86").Value = -221.0667000000001
$ws.Range("H88").Value = 7936.857
$ws.Range("I88").Value = 8750
$ws.Range("K88").Value = 8750
$ws.Range("M88").Value = -8344
$ws.Range("H89").Value = 2319.923
$ws.Range("I89").Value = 1344.0667
$ws.Range("K89").Value = 6720.333500000001
$ws.Range("M89").Value = -1104.333500000001
$ws.Range("H91").Value = 7936.857
$ws.Range("I91").Value = 8750
$ws.Range("K91").Value = 8750
$ws.Range("M91").Value = -7346
$ws.Range("H106").Value = 839.8570999999999
$ws.Range("I106").Value = 827.61536
$ws.Range("K106").Value = 827.61536
$ws.Range("M106").Value = -196.61536
$ws.Range("H107").Value = 3628.7144
$ws.Range("I107").Value = 4724.8
$ws.Range("J107").Value = 888.5
$ws.Range("K107").Value = 4724.8
$ws.Range("L107").Value = 888.5
$ws.Range("M107").Value = -2804.8
$ws.Range("N107").Value = -4728.5
$ws.Range("H112").Value = 1796.2424
$ws.Range("I112").Value = 1296.3334
$ws.Range("K112").Value = 3889.0002
$ws.Range("M112").Value = -2781.0002
$ws.Range("H113").Value = 9849.5
$ws.Range("I113").Value = 8165
$ws.Range("J113").Value = 10571.429
$ws.Range("K113").Value = 8165
$ws.Range("L113").Value = 10571.429
$ws.Range("M113").Value = -4911
$ws.Range("N113").Value = -17079.429
$ws.Range("H138").Value = 3278.9473
$ws.Range("J138").Value = 4224.143
$ws.Range("L138").Value = 12672.429
$ws.Range("N138").Value = -22952.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15390605
$ws.Range("I32").Value = 15390605
$ws.Range("K32").Value = 15390605
$ws.Range("M32").Value = -15390318
$ws.Range("H34").Value = 24828.334
$ws.Range("I34").Value = 21995
$ws.Range("K34").Value = 21995
$ws.Range("M34").Value = -21724
$ws.Range("H45").Value = 2250
$ws.Range("I45").Value = 1666.6666
$ws.Range("J45").Value = 4000
$ws.Range("K45").Value = 1666.6666
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = -1289.6666
$ws.Range("N45").Value = -4754
$ws.Range("H61").Value = 2440.25
$ws.Range("I61").Value = 2567.1
$ws.Range("J61").Value = 1806
$ws.Range("K61").Value = 2567.1
$ws.Range("L61").Value = 1806
$ws.Range("M61").Value = -2355.1
$ws.Range("N61").Value = -2230
$ws.Range("H74").Value = 2083.1035
$ws.Range("I74").Value = 1978.9286
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 1978.9286
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -1104.9286
$ws.Range("N74").Value = -6748
$ws.Range("H77").Value = 2083.1035
$ws.Range("I77").Value = 1978.9286
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 9894.643
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -5526.643
$ws.Range("N77").Value = -33736
$ws.Range("H88").Value = 1270.35
$ws.Range("J88").Value = 1300.4117
$ws.Range("L88").Value = 1300.4117
$ws.Range("N88").Value = -2112.4117
$ws.Range("H91").Value = 1270.35
$ws.Range("J91").Value = 1300.4117
$ws.Range("L91").Value = 1300.4117
$ws.Range("N91").Value = -4108.411700000001
$ws.Range("H97").Value = 2780.3572
$ws.Range("J97").Value = 6167.75
$ws.Range("L97").Value = 6167.75
$ws.Range("N97").Value = -7159.75
$ws.Range("H122").Value = 1678.7
$ws.Range("I122").Value = 1221.625
$ws.Range("J122").Value = 3507
$ws.Range("K122").Value = 3664.875
$ws.Range("L122").Value = 10521
$ws.Range("M122").Value = -1214.875
$ws.Range("N122").Value = -15421
$ws.Range("H132").Value = 2985.0588
$ws.Range("I132").Value = 2983.4
$ws.Range("J132").Value = 2997.5
$ws.Range("K132").Value = 8950.200000000001
$ws.Range("L132").Value = 8992.5
$ws.Range("M132").Value = -6420.200000000001
$ws.Range("N132").Value = -14052.5
$ws.Range("H136").Value = 2440.25
$ws.Range("I136").Value = 2567.1
$ws.Range("J136").Value = 1806
$ws.Range("K136").Value = 7701.299999999999
$ws.Range("L136").Value = 5418
$ws.Range("M136").Value = -5151.299999999999
$ws.Range("N136").Value = -10518

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3400.5
$ws.Range("I20").Value = 1208.3636
$ws.Range("J20").Value = 6845.2856
$ws.Range("K20").Value = 1208.3636
$ws.Range("L20").Value = 6845.2856
$ws.Range("M20").Value = -961.3635999999999
$ws.Range("N20").Value = -7339.2856
$ws.Range("H21").Value = 20000
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20472
$ws.Range("H63").Value = 40135.5
$ws.Range("J63").Value = 40135.5
$ws.Range("L63").Value = 40135.5
$ws.Range("N63").Value = -41507.5
$ws.Range("H66").Value = 40135.5
$ws.Range("J66").Value = 40135.5
$ws.Range("L66").Value = 120406.5
$ws.Range("N66").Value = -127270.5
$ws.Range("H86").Value = 3565.4827
$ws.Range("I86").Value = 2607.6316
$ws.Range("K86").Value = 2607.6316
$ws.Range("M86").Value = -1484.6316
$ws.Range("H89").Value = 3565.4827
$ws.Range("I89").Value = 2607.6316
$ws.Range("K89").Value = 13038.158
$ws.Range("M89").Value = -7422.158000000001
$ws.Range("H134").Value = 4200
$ws.Range("I134").Value = 3916.6667
$ws.Range("K134").Value = 11750.0001
$ws.Range("M134").Value = -9215.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1706.4667
$ws.Range("I31").Value = 1654.317
$ws.Range("J31").Value = 2241
$ws.Range("K31").Value = 1654.317
$ws.Range("L31").Value = 2241
$ws.Range("M31").Value = -1359.317
$ws.Range("N31").Value = -2831
$ws.Range("H32").Value = 4767.5
$ws.Range("I32").Value = 4767.5
$ws.Range("K32").Value = 4767.5
$ws.Range("M32").Value = -4451.5
$ws.Range("H34").Value = 1706.4667
$ws.Range("I34").Value = 1654.317
$ws.Range("J34").Value = 2241
$ws.Range("K34").Value = 1654.317
$ws.Range("L34").Value = 2241
$ws.Range("M34").Value = -1452.317
$ws.Range("N34").Value = -2645
$ws.Range("H45").Value = 9074
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 9074
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 9074
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -10260
$ws.Range("H51").Value = 38000
$ws.Range("I51").Value = 38000
$ws.Range("K51").Value = 38000
$ws.Range("M51").Value = -37264
$ws.Range("H58").Value = 1554.3903
$ws.Range("I58").Value = 870.65515
$ws.Range("J58").Value = 3206.75
$ws.Range("K58").Value = 870.65515
$ws.Range("L58").Value = 3206.75
$ws.Range("M58").Value = -667.65515
$ws.Range("N58").Value = -3612.75
$ws.Range("H61").Value = 38000
$ws.Range("I61").Value = 38000
$ws.Range("K61").Value = 38000
$ws.Range("M61").Value = -37652
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H70").Value = 28333.334
$ws.Range("J70").Value = 28333.334
$ws.Range("L70").Value = 28333.334
$ws.Range("N70").Value = -28963.334
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H73").Value = 28333.334
$ws.Range("J73").Value = 28333.334
$ws.Range("L73").Value = 28333.334
$ws.Range("N73").Value = -30517.334
$ws.Range("H74").Value = 67437.664
$ws.Range("J74").Value = 67437.664
$ws.Range("L74").Value = 67437.664
$ws.Range("N74").Value = -69185.664
$ws.Range("H77").Value = 67437.664
$ws.Range("J77").Value = 67437.664
$ws.Range("L77").Value = 202312.992
$ws.Range("N77").Value = -211048.992
$ws.Range("H115").Value = 30844.691
$ws.Range("J115").Value = 30844.691
$ws.Range("L115").Value = 30844.691
$ws.Range("N115").Value = -33194.691
$ws.Range("H121").Value = 19999
$ws.Range("J121").Value = 19999
$ws.Range("L121").Value = 19999
$ws.Range("N121").Value = -22619
$ws.Range("H132").Value = 1735.56
$ws.Range("I132").Value = 1625.6522
$ws.Range("K132").Value = 4876.9566
$ws.Range("M132").Value = -2346.9566
$ws.Range("H134").Value = 1860.5264
$ws.Range("I134").Value = 1995.4688
$ws.Range("K134").Value = 5986.4064
$ws.Range("M134").Value = -3451.4064
$ws.Range("H136").Value = 1554.3903
$ws.Range("I136").Value = 870.65515
$ws.Range("J136").Value = 3206.75
$ws.Range("K136").Value = 2611.96545
$ws.Range("L136").Value = 9620.25
$ws.Range("M136").Value = -61.96545000000015
$ws.Range("N136").Value = -14720.25
$ws.Range("H141").Value = 46350.918
$ws.Range("J141").Value = 48019.184
$ws.Range("L141").Value = 48019.184
$ws.Range("N141").Value = -58379.184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 85.333336
$ws.Range("J2").Value = 82.666664
$ws.Range("L2").Value = 495.999984
$ws.Range("N2").Value = -721.999984
$ws.Range("H7").Value = 2105
$ws.Range("I7").Value = 235
$ws.Range("J7").Value = 3975
$ws.Range("K7").Value = 705
$ws.Range("L7").Value = 11925
$ws.Range("M7").Value = -593
$ws.Range("N7").Value = -12149
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H38").Value = 1293.4375
$ws.Range("I38").Value = 26.5
$ws.Range("J38").Value = 2053.6
$ws.Range("K38").Value = 79.5
$ws.Range("L38").Value = 6160.799999999999
$ws.Range("M38").Value = 267.5
$ws.Range("N38").Value = -6854.799999999999
$ws.Range("H107").Value = 350
$ws.Range("I107").Value = 350
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1050
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 870
$ws.Range("N107").ClearContents()
$ws.Range("H109").Value = 1045.4166
$ws.Range("J109").Value = 8000
$ws.Range("L109").Value = 24000
$ws.Range("N109").Value = -26080
$ws.Range("H122").Value = 754.61536
$ws.Range("I122").Value = 941.44446
$ws.Range("J122").Value = 334.25
$ws.Range("K122").Value = 8473.00014
$ws.Range("L122").Value = 3008.25
$ws.Range("M122").Value = -6023.00014
$ws.Range("N122").Value = -7908.25
$ws.Range("H131").Value = 3719.1667
$ws.Range("I131").Value = 2414.75
$ws.Range("J131").Value = 4762.7
$ws.Range("K131").Value = 7244.25
$ws.Range("L131").Value = 14288.1
$ws.Range("M131").Value = -2204.25
$ws.Range("N131").Value = -24368.1
$ws.Range("H132").Value = 1434.2142
$ws.Range("J132").Value = 1274.9166
$ws.Range("L132").Value = 11474.2494
$ws.Range("N132").Value = -16534.2494
$ws.Range("H136").Value = 2989.2666
$ws.Range("I136").Value = 977
$ws.Range("K136").Value = 2931
$ws.Range("M136").Value = 2169
$ws.Range("H137").Value = 2733.1924
$ws.Range("I137").Value = 2343.7856
$ws.Range("K137").Value = 7031.3568
$ws.Range("M137").Value = -1931.3568
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 7143022.5
$ws.Range("I2").Value = 42.57143
$ws.Range("J2").Value = 14286002
$ws.Range("K2").Value = 42.57143
$ws.Range("L2").Value = 14286002
$ws.Range("M2").Value = 70.42857000000001
$ws.Range("N2").Value = -14286228
$ws.Range("H43").Value = 4234.3335
$ws.Range("I43").Value = 1083.2
$ws.Range("K43").Value = 1083.2
$ws.Range("M43").Value = -932.2
$ws.Range("H44").Value = 25999
$ws.Range("J44").Value = 25999
$ws.Range("L44").Value = 25999
$ws.Range("N44").Value = -27191
$ws.Range("H47").Value = 13010
$ws.Range("J47").Value = 13010
$ws.Range("L47").Value = 13010
$ws.Range("N47").Value = -14146
$ws.Range("H70").Value = 558055.5
$ws.Range("I70").Value = 558055.5
$ws.Range("K70").Value = 558055.5
$ws.Range("M70").Value = -557785.5
$ws.Range("H73").Value = 558055.5
$ws.Range("I73").Value = 558055.5
$ws.Range("K73").Value = 558055.5
$ws.Range("M73").Value = -557119.5
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31748
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -98736
$ws.Range("H92").Value = 18243.428
$ws.Range("J92").Value = 18243.428
$ws.Range("L92").Value = 18243.428
$ws.Range("N92").Value = -21987.428
$ws.Range("H122").Value = 5983.5264
$ws.Range("I122").Value = 6212.5
$ws.Range("J122").Value = 5729.1113
$ws.Range("K122").Value = 18637.5
$ws.Range("L122").Value = 17187.3339
$ws.Range("M122").Value = -16187.5
$ws.Range("N122").Value = -22087.3339
$ws.Range("H123").Value = 33785.715
$ws.Range("J123").Value = 33785.715
$ws.Range("L123").Value = 33785.715
$ws.Range("N123").Value = -38685.715
$ws.Range("H126").Value = 9933.333000000001
$ws.Range("I126").Value = 9950
$ws.Range("K126").Value = 29850
$ws.Range("M126").Value = -27380
$ws.Range("H132").Value = 2589.4
$ws.Range("I132").Value = 2321.9092
$ws.Range("J132").Value = 3325
$ws.Range("K132").Value = 6965.7276
$ws.Range("L132").Value = 9975
$ws.Range("M132").Value = -4435.7276
$ws.Range("N132").Value = -15035

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("M16").Value = -830
$ws.Range("H40").Value = 9045.6875
$ws.Range("J40").Value = 4686.3335
$ws.Range("L40").Value = 4686.3335
$ws.Range("N40").Value = -4958.3335
$ws.Range("H46").Value = 4219.625
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 4219.625
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 4219.625
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -4595.625
$ws.Range("H61").Value = 1548.9286
$ws.Range("I61").Value = 1257.2609
$ws.Range("J61").Value = 2890.6
$ws.Range("K61").Value = 1257.2609
$ws.Range("L61").Value = 2890.6
$ws.Range("M61").Value = -1055.2609
$ws.Range("N61").Value = -3294.6
$ws.Range("H113").Value = 1548.9286
$ws.Range("I113").Value = 1257.2609
$ws.Range("J113").Value = 2890.6
$ws.Range("K113").Value = 1257.2609
$ws.Range("L113").Value = 2890.6
$ws.Range("M113").Value = 912.7391
$ws.Range("N113").Value = -7230.6
$ws.Range("H132").Value = 5015.84
$ws.Range("I132").Value = 3254.7778
$ws.Range("K132").Value = 9764.3334
$ws.Range("M132").Value = -7234.3334
$ws.Range("H136").Value = 4505.609
$ws.Range("I136").Value = 4515.778
$ws.Range("K136").Value = 13547.334
$ws.Range("M136").Value = -10997.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 25018738
$ws.Range("I2").Value = 50025000
$ws.Range("J2").Value = 12474.5
$ws.Range("K2").Value = 50025000
$ws.Range("L2").Value = 12474.5
$ws.Range("M2").Value = -50024888
$ws.Range("N2").Value = -12698.5
$ws.Range("H62").Value = 3499.3333
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 3499
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 3499
$ws.Range("M62").Value = -2876
$ws.Range("N62").Value = -4747
$ws.Range("H65").Value = 3499.3333
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 3499
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 17495
$ws.Range("M65").Value = -14380
$ws.Range("N65").Value = -23735
$ws.Range("H132").Value = 2099
$ws.Range("I132").Value = 2124.3572
$ws.Range("J132").Value = 2010.25
$ws.Range("K132").Value = 6373.071599999999
$ws.Range("L132").Value = 6030.75
$ws.Range("M132").Value = -3843.071599999999
$ws.Range("N132").Value = -11090.75
